# The sheet originally contained a small header row (r=2: "municipios"/"casos"/"óbitos")
# directly above the municipality data, and two extra aggregate rows at the very
# bottom (r=89: "outros estados", r=90: "outros paises"). The update removes that
# header row (so the municipality data shifts up to start right under the
# "Unnamed: n" row) and drops the two trailing aggregate rows, shrinking the used
# range from A1:C90 down to A1:C87.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the bottom two rows first (order doesn't matter since they are below
# row 2, but doing the bottom rows first keeps row numbers for them stable).
$ws.Range("A89:C90").EntireRow.Delete() | Out-Null

# Remove the redundant header row right above the data (was row 2: municipios/casos/óbitos).
$ws.Range("A2:C2").EntireRow.Delete() | Out-Null
